$p = $ppt.ActivePresentation

# Slide 1 ("Nodejs Basics" title slide) - the subtitle placeholder's third
# paragraph reads "04 – Node Modules". Split it into two runs so the
# session title becomes "04 – Express": keep "04 " in the original run,
# and replace the remainder ("– Node Modules") with "– Express" - this is
# the portion a presenter would retype when moving on to the new session.
$slide = $p.Slides.Item(1)
$subtitle = $slide.Shapes.Item(3)
$tr = $subtitle.TextFrame.TextRange
$sessionPara = $tr.Paragraphs(3, 1)

# Characters(4, 14) selects "– Node Modules" (everything after "04 ").
$tail = $sessionPara.Characters(4, 14)
$tail.Text = "– Express"
